$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1 & 2: move the "_GoBack" bookmark from right after the title
#        ("Space Escape Mechanics and Meaningful Play") down to right
#        after the " a menu selection. " run (end of that paragraph).
#
# The engine's Bookmarks.Add() mis-handles a zero-length range that
# sits exactly at the end-of-paragraph-text boundary (right before the
# paragraph mark), so we work around it: temporarily insert a marker
# character at that boundary, wrap the bookmark around it, then delete
# the marker again. Re-using the existing bookmark name ("_GoBack")
# relocates it, so the old bookmarkStart/bookmarkEnd pair near the
# title disappears automatically.
# ------------------------------------------------------------------

$anchor = $d.Content
[void]$anchor.Find.Execute(" a menu selection. ")
$anchor.Collapse(0)
$pos = $anchor.Start

$markerChar = [char]1
$marker = $d.Range($pos, $pos)
$marker.InsertAfter($markerChar)

$bmRange = $d.Range($pos, $pos + 1)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($pos, $pos + 1).Text = ""

# ------------------------------------------------------------------
# 3: remove the two <w:proofErr/> markers around "at" and merge the
#    three runs back into a single run (re-typing the same text over
#    the found range collapses it to one run with uniform formatting).
# ------------------------------------------------------------------

$old = "It is more of a distance measurement that is slightly better at evaluating the user’s challenge performance"
$new = "It is more of a distance measurement that is slightly better at evaluating the user’s challenge performance"
[void]$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
